$wb = $excel.ActiveWorkbook

# ALC row 19
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 483
$ws.Range("J19").Value = 474.75
$ws.Range("L19").Value = 474.75
$ws.Range("N19").Value = -824.75

# ALC row 51
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 10999.5
$ws.Range("I51").Value = 10999
$ws.Range("K51").Value = 10999
$ws.Range("M51").Value = -10515

# ALC row 80
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 1879593.2
$ws.Range("I80").Value = 1433.6666
$ws.Range("J80").Value = 2201563.5
$ws.Range("K80").Value = 4300.9998
$ws.Range("L80").Value = 6604690.5
$ws.Range("M80").Value = -3302.9998
$ws.Range("N80").Value = -6606686.5

# ALC row 83
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 1879593.2
$ws.Range("I83").Value = 1433.6666
$ws.Range("J83").Value = 2201563.5
$ws.Range("K83").Value = 12902.9994
$ws.Range("L83").Value = 19814071.5
$ws.Range("M83").Value = -7910.999400000001
$ws.Range("N83").Value = -19824055.5

# ALC row 125
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 661.6667
$ws.Range("I125").Value = 223.33333
$ws.Range("J125").Value = 1100
$ws.Range("K125").Value = 2009.99997
$ws.Range("L125").Value = 9900
$ws.Range("M125").Value = 450.0000300000002
$ws.Range("N125").Value = -14820

# ALC row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 790.8570999999999
$ws.Range("I129").Value = 647.125
$ws.Range("J129").Value = 848.35
$ws.Range("K129").Value = 1941.375
$ws.Range("L129").Value = 2545.05
$ws.Range("M129").Value = 3058.625
$ws.Range("N129").Value = -12545.05

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 69645.92999999999
$ws.Range("I137").Value = 3699.7144
$ws.Range("J137").Value = 127348.875
$ws.Range("K137").Value = 11099.1432
$ws.Range("L137").Value = 382046.625
$ws.Range("M137").Value = -8549.143199999999
$ws.Range("N137").Value = -387146.625

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5056.125
$ws.Range("I61").Value = 2224.5
$ws.Range("J61").Value = 6000
$ws.Range("K61").Value = 2224.5
$ws.Range("L61").Value = 6000
$ws.Range("M61").Value = -2012.5
$ws.Range("N61").Value = -6424

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 5056.125
$ws.Range("I136").Value = 2224.5
$ws.Range("J136").Value = 6000
$ws.Range("K136").Value = 6673.5
$ws.Range("L136").Value = 18000
$ws.Range("M136").Value = -4123.5
$ws.Range("N136").Value = -23100

# BSM row 64
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 45455220
$ws.Range("I64").Value = 111112216
$ws.Range("J64").Value = 376
$ws.Range("K64").Value = 111112216
$ws.Range("L64").Value = 376
$ws.Range("M64").Value = -111111991
$ws.Range("N64").Value = -826

# BSM row 67
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H67").Value = 45455220
$ws.Range("I67").Value = 111112216
$ws.Range("J67").Value = 376
$ws.Range("K67").Value = 111112216
$ws.Range("L67").Value = 376
$ws.Range("M67").Value = -111111436
$ws.Range("N67").Value = -1936

# BSM row 132
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 47500
$ws.Range("J132").Value = 47500
$ws.Range("L132").Value = 47500
$ws.Range("N132").Value = -57620

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 97246.55
$ws.Range("I134").Value = 97246.55
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 291739.65
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -289204.65
$ws.Range("N134").ClearContents() | Out-Null

# CRP row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1246.5714
$ws.Range("I16").Value = 1056.5
$ws.Range("K16").Value = 1056.5
$ws.Range("M16").Value = -769.5

# CRP row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1246.5714
$ws.Range("I113").Value = 1056.5
$ws.Range("K113").Value = 1056.5
$ws.Range("M113").Value = 1113.5

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 26530.191
$ws.Range("I132").Value = 36942.5
$ws.Range("J132").Value = 5705.5713
$ws.Range("K132").Value = 110827.5
$ws.Range("L132").Value = 17116.7139
$ws.Range("M132").Value = -108297.5
$ws.Range("N132").Value = -22176.7139

# CUL row 12
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 74.07143000000001
$ws.Range("I12").Value = 20
$ws.Range("K12").Value = 60
$ws.Range("M12").Value = 113

# CUL row 39
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 2292.5454
$ws.Range("J39").Value = 2479.7778
$ws.Range("L39").Value = 7439.3334
$ws.Range("N39").Value = -8027.3334

# CUL row 98
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 540.875
$ws.Range("I98").Value = 382.5
$ws.Range("K98").Value = 1147.5
$ws.Range("M98").Value = 350.5

# CUL row 121
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 5870.421
$ws.Range("I121").Value = 467
$ws.Range("J121").Value = 7800.2144
$ws.Range("K121").Value = 1401
$ws.Range("L121").Value = 23400.6432
$ws.Range("M121").Value = -91
$ws.Range("N121").Value = -26020.6432

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 779.02
$ws.Range("J131").Value = 784.40814
$ws.Range("L131").Value = 2353.22442
$ws.Range("N131").Value = -12433.22442

# GSM row 15
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents() | Out-Null

# GSM row 81
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents() | Out-Null

# GSM row 84
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents() | Out-Null

# GSM row 107
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1388.2
$ws.Range("I107").Value = 277.6
$ws.Range("J107").Value = 2498.8
$ws.Range("K107").Value = 277.6
$ws.Range("L107").Value = 2498.8
$ws.Range("M107").Value = 1642.4
$ws.Range("N107").Value = -6338.8

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 111257.86
$ws.Range("I132").Value = 115956.89
$ws.Range("J132").Value = 102799.6
$ws.Range("K132").Value = 347870.67
$ws.Range("L132").Value = 308398.8
$ws.Range("M132").Value = -345340.67
$ws.Range("N132").Value = -313458.8

# LTW row 59
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H59").Value = 29700
$ws.Range("J59").Value = 29700
$ws.Range("L59").Value = 29700
$ws.Range("N59").Value = -31008

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 1229158.2
$ws.Range("I122").Value = 2453980
$ws.Range("K122").Value = 7361940
$ws.Range("M122").Value = -7359490

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 35806.535
$ws.Range("I136").Value = 47599.816
$ws.Range("J136").Value = 3375
$ws.Range("K136").Value = 142799.448
$ws.Range("L136").Value = 10125
$ws.Range("M136").Value = -140249.448
$ws.Range("N136").Value = -15225

# WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1093.7084
$ws.Range("J126").Value = 2083.5715
$ws.Range("L126").Value = 6250.7145
$ws.Range("N126").Value = -11190.7145

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 922774.5
$ws.Range("I136").Value = 1241481
$ws.Range("J136").Value = 2067
$ws.Range("K136").Value = 3724443
$ws.Range("L136").Value = 6201
$ws.Range("M136").Value = -3721893
$ws.Range("N136").Value = -11301

Write-Host "Applied 31 row updates across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets"
